$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 51
$ws.Range("H51").Value = 15934.786
$ws.Range("I51").Value = 100001
$ws.Range("J51").Value = 1923.75
$ws.Range("K51").Value = 100001
$ws.Range("L51").Value = 1923.75
$ws.Range("M51").Value = -99517
$ws.Range("N51").Value = -2891.75

# Row 112
$ws.Range("H112").Value = 16112.642
$ws.Range("J112").Value = 16725.883
$ws.Range("L112").Value = 50177.649
$ws.Range("N112").Value = -52393.649

# Row 137
$ws.Range("H137").Value = 1056.1343
$ws.Range("I137").Value = 917
$ws.Range("J137").Value = 1582.8572
$ws.Range("K137").Value = 2751
$ws.Range("L137").Value = 4748.571599999999
$ws.Range("M137").Value = -201
$ws.Range("N137").Value = -9848.571599999999

# Row 138
$ws.Range("H138").Value = 2085.4893
$ws.Range("I138").Value = 1316.1724
$ws.Range("J138").Value = 3324.9443
$ws.Range("K138").Value = 3948.5172
$ws.Range("L138").Value = 9974.832900000001
$ws.Range("M138").Value = 1191.4828
$ws.Range("N138").Value = -20254.8329

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 982.9322
$ws.Range("I61").Value = 910.38635
$ws.Range("J61").Value = 1195.7333
$ws.Range("K61").Value = 910.38635
$ws.Range("L61").Value = 1195.7333
$ws.Range("M61").Value = -698.38635
$ws.Range("N61").Value = -1619.7333

# Row 74
$ws.Range("H74").Value = 1150.7142
$ws.Range("I74").Value = 1337.4814
$ws.Range("J74").Value = 814.5333000000001
$ws.Range("K74").Value = 1337.4814
$ws.Range("L74").Value = 814.5333000000001
$ws.Range("M74").Value = -463.4813999999999
$ws.Range("N74").Value = -2562.5333

# Row 77
$ws.Range("H77").Value = 1150.7142
$ws.Range("I77").Value = 1337.4814
$ws.Range("J77").Value = 814.5333000000001
$ws.Range("K77").Value = 6687.406999999999
$ws.Range("L77").Value = 4072.6665
$ws.Range("M77").Value = -2319.406999999999
$ws.Range("N77").Value = -12808.6665

# Row 136
$ws.Range("H136").Value = 982.9322
$ws.Range("I136").Value = 910.38635
$ws.Range("J136").Value = 1195.7333
$ws.Range("K136").Value = 2731.15905
$ws.Range("L136").Value = 3587.199900000001
$ws.Range("M136").Value = -181.1590500000002
$ws.Range("N136").Value = -8687.1999

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5297281.5
$ws.Range("I134").Value = 1537
$ws.Range("J134").Value = 10111595
$ws.Range("K134").Value = 4611
$ws.Range("L134").Value = 30334785
$ws.Range("M134").Value = -2076
$ws.Range("N134").Value = -30339855

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1160.7701
$ws.Range("I31").Value = 904.35596
$ws.Range("J31").Value = 1701.0714
$ws.Range("K31").Value = 904.35596
$ws.Range("L31").Value = 1701.0714
$ws.Range("M31").Value = -609.35596
$ws.Range("N31").Value = -2291.0714

# Row 34
$ws.Range("H34").Value = 1160.7701
$ws.Range("I34").Value = 904.35596
$ws.Range("J34").Value = 1701.0714
$ws.Range("K34").Value = 904.35596
$ws.Range("L34").Value = 1701.0714
$ws.Range("M34").Value = -702.35596
$ws.Range("N34").Value = -2105.0714

# Row 58
$ws.Range("H58").Value = 18519826
$ws.Range("I58").Value = 24391714
$ws.Range("J58").Value = 793.38464
$ws.Range("K58").Value = 24391714
$ws.Range("L58").Value = 793.38464
$ws.Range("M58").Value = -24391511
$ws.Range("N58").Value = -1199.38464

# Row 132
$ws.Range("H132").Value = 16668336
$ws.Range("I132").Value = 1531.7333
$ws.Range("K132").Value = 4595.199900000001
$ws.Range("M132").Value = -2065.199900000001

# Row 133
$ws.Range("H133").Value = 30326
$ws.Range("J133").Value = 30326
$ws.Range("L133").Value = 30326
$ws.Range("N133").Value = -35386

# Row 134
$ws.Range("H134").Value = 1193.8857
$ws.Range("I134").Value = 1218
$ws.Range("K134").Value = 3654
$ws.Range("M134").Value = -1119

# Row 136
$ws.Range("H136").Value = 18519826
$ws.Range("I136").Value = 24391714
$ws.Range("J136").Value = 793.38464
$ws.Range("K136").Value = 73175142
$ws.Range("L136").Value = 2380.15392
$ws.Range("M136").Value = -73172592
$ws.Range("N136").Value = -7480.15392

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 3969.6667
$ws.Range("J3").Value = 7000
$ws.Range("L3").Value = 21000
$ws.Range("N3").Value = -21224

# Row 12
$ws.Range("H12").Value = 24409.83
$ws.Range("I12").Value = 3.3333333
$ws.Range("J12").Value = 31274.156
$ws.Range("K12").Value = 9.999999900000001
$ws.Range("L12").Value = 93822.46799999999
$ws.Range("M12").Value = 163.0000001
$ws.Range("N12").Value = -94168.46799999999

# Row 131
$ws.Range("H131").Value = 786.2929
$ws.Range("J131").Value = 808.8461
$ws.Range("L131").Value = 2426.5383
$ws.Range("N131").Value = -12506.5383

# Row 132
$ws.Range("H132").Value = 13161584
$ws.Range("I132").Value = 656.1667
$ws.Range("J132").Value = 25006420
$ws.Range("K132").Value = 5905.5003
$ws.Range("L132").Value = 225057780
$ws.Range("M132").Value = -3375.5003
$ws.Range("N132").Value = -225062840

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 4195.44
$ws.Range("I55").Value = 8544.583000000001
$ws.Range("J55").Value = 180.84616
$ws.Range("K55").Value = 8544.583000000001
$ws.Range("L55").Value = 180.84616
$ws.Range("M55").Value = -8371.583000000001
$ws.Range("N55").Value = -526.8461600000001

# Row 132
$ws.Range("H132").Value = 29420168
$ws.Range("I132").Value = 55558292
$ws.Range("J132").Value = 14775.3125
$ws.Range("K132").Value = 166674876
$ws.Range("L132").Value = 44325.9375
$ws.Range("M132").Value = -166672346
$ws.Range("N132").Value = -49385.9375

# Row 136
$ws.Range("H136").Value = 1982.909
$ws.Range("I136").Value = 2132.75
$ws.Range("J136").Value = 1583.3334
$ws.Range("K136").Value = 6398.25
$ws.Range("L136").Value = 4750.0002
$ws.Range("M136").Value = -3848.25
$ws.Range("N136").Value = -9850.0002

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 282.14285
$ws.Range("I107").Value = 227.625
$ws.Range("J107").Value = 456.6
$ws.Range("K107").Value = 682.875
$ws.Range("L107").Value = 1369.8
$ws.Range("M107").Value = 1237.125
$ws.Range("N107").Value = -5209.8

# Row 132
$ws.Range("H132").Value = 26932.979
$ws.Range("I132").Value = 39365.52
$ws.Range("J132").Value = 8284.166999999999
$ws.Range("K132").Value = 118096.56
$ws.Range("L132").Value = 24852.501
$ws.Range("M132").Value = -115566.56
$ws.Range("N132").Value = -29912.501

# Row 136
$ws.Range("H136").Value = 7466341
$ws.Range("I136").Value = 12200138
$ws.Range("K136").Value = 36600414
$ws.Range("M136").Value = -36597864
